$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates -------------------------------------------------------

# D2: investigator name with qualification changes from the DeLuca entry
# to the Bansal entry. It also picks up a smaller (9pt) grey Arial font,
# matching the one already used on the "Medical License Number" cell (P2),
# just at 9pt instead of 11pt.
$ws.Range("P2").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("D2").Font.Size = 9
$ws.Range("D2").Value = "Bansal, Padam C. MD"

# G2: Last Name DeLuca -> Padam (first name of the new investigator)
$ws.Range("G2").Value = "Padam"

# I2: First Name William -> Bansal (last name of the new investigator)
$ws.Range("I2").Value = "Bansal"

# --- New rows for additional (Sub) investigators --------------------------

$ws.Range("A3").Value = "Sub"
$ws.Range("D3").Value = "fogari, robert MD"
$ws.Range("G3").Value = "Robert"
$ws.Range("I3").Value = "Fogari"

$ws.Range("A4").Value = "Sub"
$ws.Range("D4").Value = "James W. Michalek"
$ws.Range("G4").Value = "James"
$ws.Range("I4").Value = "Michalek"

$ws.Range("A5").Value = "Sub"
$ws.Range("D5").Value = "Kincaid, William Ralph"
$ws.Range("G5").Value = "William"
$ws.Range("H5").Value = "Kincaid"
$ws.Range("I5").Value = "Ralph"

# C2: was blank -> "0102" entered as text (leading apostrophe keeps the
# leading zero and marks the cell as quote-prefixed text).
$ws.Range("C2").Value = "'0102"

# --- Selection / view ------------------------------------------------------
$ws.Range("C2").Select() | Out-Null
